$wb = $excel.ActiveWorkbook

# --- Update "expense_bills" sheet: prefix file names with "email_" ---
$wsBills = $wb.Worksheets.Item("expense_bills")
for ($row = 2; $row -le 6; $row++) {
    $cell = $wsBills.Cells.Item($row, 1)
    $oldValue = $cell.Value()
    $cell.Value = "email_" + $oldValue
}

# --- Append two new rows to "unmatched_expenses" sheet ---
$wsUnmatched = $wb.Worksheets.Item("unmatched_expenses")

$wsUnmatched.Cells.Item(3, 1).Value = "UNKNOWN DATE"
$wsUnmatched.Cells.Item(3, 2).Value = 780
$wsUnmatched.Cells.Item(3, 3).Value = "Amount in Words: 780 Rupees Only (parsing failed)"

$wsUnmatched.Cells.Item(4, 1).Value = "UNKNOWN DATE"
$wsUnmatched.Cells.Item(4, 2).Value = 780
$wsUnmatched.Cells.Item(4, 3).Value = "Amount in Words: 780 Rupees Only (parsing failed)"
